$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap full row contents (columns B..AB) between paired rows.
# (Column A is a fixed sequential index and is left untouched.)
$cols = 2..28
$rowPairs = @(
    @(6, 7),
    @(37, 38),
    @(41, 42),
    @(126, 127),
    @(138, 139),
    @(142, 143),
    @(148, 150)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    foreach ($c in $cols) {
        $cell1 = $ws.Cells.Item($r1, $c)
        $cell2 = $ws.Cells.Item($r2, $c)
        $v1 = $cell1.Value()
        $v2 = $cell2.Value()
        $cell1.Value = $v2
        $cell2.Value = $v1
    }
}
